$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = 9

$ws.Cells.Item($r, 1).Value = 45862.72465423148
$ws.Cells.Item($r, 1).NumberFormat = $ws.Cells.Item(5, 1).NumberFormat

$ws.Cells.Item($r, 2).Value = 2025
$ws.Cells.Item($r, 3).Value = 30
$ws.Cells.Item($r, 4).Value = 19.33
$ws.Cells.Item($r, 5).Value = 74.44
$ws.Cells.Item($r, 6).Value = 87.89
$ws.Cells.Item($r, 7).Value = 13.54
$ws.Cells.Item($r, 8).Value = "ESE"
$ws.Cells.Item($r, 9).Value = 0
$ws.Cells.Item($r, 10).Value = "17:23:30"
